# Add participation flag seeding based on existing deferral rates.
# Re-seeded the Participating headcount (and downstream aggregates) for
# several projection years, which changes Participating, Participation
# Rate, Avg Deferral Rate, contribution totals, compensation totals and
# the resulting cost-percentage ratios.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Year 1) ---------------------------------------------------
$ws.Range("G2").Value = 0.1064597248414084
$ws.Range("H2").Value = 0.08976016016040311
$ws.Range("I2").Value = 523152.5604975562
$ws.Range("J2").Value = 198962.2804497781
$ws.Range("L2").Value = 198962.2804497781
$ws.Range("M2").Value = 722114.8409473342
$ws.Range("N2").Value = 10195939.0888
$ws.Range("O2").Value = 9788198.158699997
$ws.Range("P2").Value = 0.01951387495717128
$ws.Range("Q2").Value = 0.02032675240365208

# --- Row 3 (Year 2) ---------------------------------------------------
$ws.Range("D3").Value = 89
$ws.Range("E3").Value = 0.8640776699029126
$ws.Range("F3").Value = 0.8640776699029126
$ws.Range("G3").Value = 0.1097230768598265
$ws.Range("H3").Value = 0.09480926058761706
$ws.Range("I3").Value = 597470.8937858229
$ws.Range("J3").Value = 233523.0270999415
$ws.Range("L3").Value = 233523.0270999415
$ws.Range("M3").Value = 830993.9208857642
$ws.Range("N3").Value = 10483919.652964
$ws.Range("O3").Value = 10076546.494961
$ws.Range("P3").Value = 0.02227440068504532
$ws.Range("Q3").Value = 0.02317490692041364

# --- Row 4 (Year 3) ---------------------------------------------------
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 0.8461538461538461
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.1143843363215895
$ws.Range("H4").Value = 0.09678674611826801
$ws.Range("I4").Value = 662421.2148822283
$ws.Range("J4").Value = 258911.8150543551
$ws.Range("L4").Value = 258911.8150543551
$ws.Range("M4").Value = 921333.0299365834
$ws.Range("N4").Value = 10944178.12365292
$ws.Range("O4").Value = 10535833.77090983
$ws.Range("P4").Value = 0.02365749279014258
$ws.Range("Q4").Value = 0.02457440205342159

# --- Row 5 (Year 4) ---------------------------------------------------
$ws.Range("C5").Value = 105
$ws.Range("D5").Value = 90
$ws.Range("E5").Value = 0.8571428571428571
$ws.Range("F5").Value = 0.8571428571428571
$ws.Range("G5").Value = 0.1137241401623681
$ws.Range("H5").Value = 0.09747783442488696
$ws.Range("I5").Value = 688290.7332236361
$ws.Range("J5").Value = 269272.6104534561
$ws.Range("L5").Value = 269272.6104534561
$ws.Range("M5").Value = 957563.3436770922
$ws.Range("N5").Value = 11180172.36496251
$ws.Range("O5").Value = 10769477.68163713
$ws.Range("P5").Value = 0.0240848353373628
$ws.Range("Q5").Value = 0.02500331199094166

# --- Row 6 (Year 5) ---------------------------------------------------
$ws.Range("D6").Value = 92
$ws.Range("E6").Value = 0.8679245283018868
$ws.Range("F6").Value = 0.8679245283018868
$ws.Range("G6").Value = 0.111601472672212
$ws.Range("H6").Value = 0.09686165552682552
$ws.Range("I6").Value = 716003.1770682526
$ws.Range("J6").Value = 280067.6496910135
$ws.Range("L6").Value = 280067.6496910135
$ws.Range("M6").Value = 996070.826759266
$ws.Range("N6").Value = 11684749.18171138
$ws.Range("O6").Value = 11270283.65788624
$ws.Range("P6").Value = 0.02396864882041002
$ws.Range("Q6").Value = 0.02457440205342159
